$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: fill in the remaining data for the "55/25/3/2.5" row (A42 already set) ---
$ws.Range("C42").Value = 103663134
$ws.Range("D42").Value = 103862
$ws.Range("E42").Value = 876998
$ws.Range("F42").Value = 103161
$ws.Range("G42").Value = 327360

# --- Extend the H-column check-sum formula (shared) down through the new rows 42:43 ---
$ws.Range("H35:H43").Formula = "=C35+10*D35+50*E35+100*F35+200*G35"

# Re-apply the existing highlight (yellow fill, same as rows 40/42 originals) to the
# appropriate cells in rows 42-43
$ws.Range("E42:H42").Interior.Color = 65535

# --- Row 43: a brand-new summary row "55/30/3/2" ---
$ws.Range("A43").Value = "55/30/3/2"
$ws.Range("A43").Font.Color = 255

$ws.Range("C43").Value = 102516230
$ws.Range("D43").Value = 214745
$ws.Range("E43").Value = 974816
$ws.Range("F43").Value = 88017
$ws.Range("G43").Value = 327360

$ws.Range("C43:D43").Interior.Color = 65535
$ws.Range("G43:H43").Interior.Color = 65535

$ws.Range("I43").Formula = "=H43/`$H`$35"
$ws.Range("I43").Interior.Color = 65535

# --- Row 46: threshold labels above the second table ---
$ws.Range("E46").Value = "100-N"
$ws.Range("F46").Value = "800-N"

# --- Row 47: header row for the second table ---
$ws.Range("A47").Value = "real"
$ws.Range("B47").Value = "sql"
$ws.Range("C47").Value = "RF2"
$ws.Range("D47").Value = "RF1"
$ws.Range("E47").Value = "RF3"
$ws.Range("F47").Value = "RF4"
$ws.Range("G47").Value = "NB"
$ws.Range("H47").Value = "LR"

# --- Rows 48-52: raw counted data for the second table ---
$ws.Range("A48").Value = 21385661
$ws.Range("B48").Value = 20502785
$ws.Range("C48").Value = 20724442
$ws.Range("D48").Value = 21340331
$ws.Range("E48").Value = 20618296
$ws.Range("F48").Value = 20619651
$ws.Range("G48").Value = 7837394
$ws.Range("H48").Value = 20710004

$ws.Range("A49").Value = 289322
$ws.Range("B49").Value = 42782
$ws.Range("C49").Value = 1
$ws.Range("D49").Value = 16
$ws.Range("E49").Value = 15
$ws.Range("F49").Value = 1
$ws.Range("G49").Value = 257785

$ws.Range("A50").Value = 299938
$ws.Range("B50").Value = 195288
$ws.Range("C50").Value = 249681
$ws.Range("D50").Value = 133074
$ws.Range("E50").Value = 253541
$ws.Range("F50").Value = 253062
$ws.Range("G50").Value = 100536
$ws.Range("H50").Value = 205930

$ws.Range("A51").Value = 60493
$ws.Range("B51").Value = 17726
$ws.Range("C51").Value = 83
$ws.Range("D51").Value = 647
$ws.Range("E51").Value = 1528
$ws.Range("F51").Value = 17
$ws.Range("G51").Value = 7609

$ws.Range("A52").Value = 74789
$ws.Range("B52").Value = 65890
$ws.Range("C52").Value = 62220
$ws.Range("D52").Value = 43875
$ws.Range("E52").Value = 65820
$ws.Range("F52").Value = 66566
$ws.Range("G52").Value = 30482
$ws.Range("H52").Value = 64796

# --- Row 53: weighted checksum per column (individual formulas, not shared) ---
$ws.Range("A53").Formula = "=A48+A49*10+A50*50+A51*100+A52*200"
$ws.Range("B53").Formula = "=B48+B49*10+B50*50+B51*100+B52*200"
$ws.Range("C53").Formula = "=C48+C49*10+C50*50+C51*100+C52*200"
$ws.Range("D53").Formula = "=D48+D49*10+D50*50+D51*100+D52*200"
$ws.Range("E53").Formula = "=E48+E49*10+E50*50+E51*100+E52*200"
$ws.Range("F53").Formula = "=F48+F49*10+F50*50+F51*100+F52*200"
$ws.Range("G53").Formula = "=G48+G49*10+G50*50+G51*100+G52*200"
$ws.Range("H53").Formula = "=H48+H49*10+H50*50+H51*100+H52*200"

# --- Row 54: each column's share of the A53 total ---
$ws.Range("B54").Formula = "=B53/`$A`$53"
$ws.Range("C54").Formula = "=C53/`$A`$53"
$ws.Range("E54").Formula = "=E53/`$A`$53"
$ws.Range("F54").Formula = "=F53/`$A`$53"
$ws.Range("G54").Formula = "=G53/`$A`$53"
$ws.Range("H54").Formula = "=H53/`$A`$53"

# --- Window / selection state ---
$ws.Range("H49").Select()
